$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Datos" to "Sheet"
$ws.Name = "Sheet"

# Re-apply the (General) number format on the creation-date placeholder cell (A2).
$ws.Range("A2").NumberFormat = "General"

# Switch the custom "dd/mm/yyyy" date format used on the date-column placeholder
# cell (B7) over to Excel's built-in short-date format.
$ws.Range("B7").NumberFormat = "mm-dd-yy"
